# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: update status text, clear the stale
# "handback not latest" error, bump the Latest Handback DateTime stamps,
# and let the Status column auto-widen for the longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet row 2 ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-06 19:07:32"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet row 2 ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-06 19:07:41"
$wsDeDe.Range("P2").Value = ""

# --- Column widths: Status columns widen for the longer text, the
#     now-empty Error Detail columns shrink back down ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.144371396019366
$wsOverview.Columns.Item(6).ColumnWidth = 29.144371396019366

$wsZhCn.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsZhCn.Columns.Item(16).ColumnWidth = 12.913719540550566

$wsDeDe.Columns.Item(3).ColumnWidth = 29.144371396019366
$wsDeDe.Columns.Item(16).ColumnWidth = 12.913719540550566
